$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: Jitter camera straffing bug - mark Resolved, update Solution text
# (written first so its new shared-string lands at the lower index, matching
# the target shared-string table order)
$ws.Range("F3").Value = "Camera overhaul, utilising Cinemachine"
$ws.Range("E3").Value = "Resolved"

# Row 2: Camera Limbo bug - mark Resolved, update Solution text
$ws.Range("F2").Value = "Camera overhaul. Clamp before value is set"
$ws.Range("E2").Value = "Resolved"

# Update the active selection to E2 (no frozen/scrolled top-left cell)
$ws.Range("E2").Select()
